$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.781.40'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.701.64'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.03'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3956'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4092'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.507'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.004'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.07'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08922'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.736'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +6.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.32'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.196'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001333'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.712.40'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '99.74'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07141'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.01'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.249'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.007'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.763.25'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.116'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.340'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.11'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.246'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +23.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '165.21'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '139.26'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.190'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.077'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +12.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09124'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +6.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.082'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03087'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +12.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2820'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.14'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.83%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.69%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.52'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09292'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.44%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.478'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7805'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.13'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.649'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7264'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.243'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.370'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.56%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.09'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '92.97'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08022'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.37%  '
